$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Update the two "published date" header labels in row 9 (G9, H9)
$ws.Range("G9").Value = "1402-03-13 (10)"
$ws.Range("H9").Value = "1402-03-13 (2)"

# Update the yearly figures in column H (most recent period)
$ws.Range("H14").Value = -26784
$ws.Range("H17").Value = 36510
$ws.Range("H18").Value = -22161
$ws.Range("H20").Value = 91996
$ws.Range("H21").Value = -2104
$ws.Range("H22").Value = 89893
$ws.Range("H24").Value = 89893
